# Rempel_2018_4-9.xlsx verification update - path issue fixed
# Adds the two new result columns (J, K) with sample data in rows 8-11,
# gives them a dedicated 6-decimal number format, sizes the columns,
# updates the current selection/view, and sets the page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data in columns J (10) and K (11), rows 8-11 ---------------------
$rng = $ws.Range("J8:K11")
$rng.Locked = $false
$rng.NumberFormat = "0.000000"

$ws.Cells.Item(8, 10).Value  = 0
$ws.Cells.Item(8, 11).Value  = 0
$ws.Cells.Item(9, 10).Value  = 0.26168224299065418
$ws.Cells.Item(9, 11).Value  = 10.3125
$ws.Cells.Item(10, 10).Value = 3.3457943925233646
$ws.Cells.Item(10, 11).Value = 14.625
$ws.Cells.Item(11, 10).Value = 0.22429906542056074
$ws.Cells.Item(11, 11).Value = 4.8125

# --- Column widths for the new columns -------------------------------------
$ws.Columns.Item(10).ColumnWidth = 8.5
$ws.Columns.Item(11).ColumnWidth = 9.5

# --- View / selection update -------------------------------------------------
# Scroll so column E is the left-most visible column, then select J13 (best
# effort - some hosts don't persist the scroll position, only the selection).
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J13").Select() | Out-Null

# Reposition the application window (best effort - window chrome placement is
# host/session specific and may not round-trip through every COM host).
$excel.Left = 6520
$excel.Top = 1020

# --- Page setup for printing ------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
